$d = $word.ActiveDocument

function Get-ParaText($index) {
    $t = $d.Paragraphs.Item($index).Range.Text
    # Strip trailing paragraph-mark character
    return $t.Substring(0, $t.Length - 1)
}

# Returns the text of the run-range that sits between two labels inside
# paragraph #17 ("Avaliação" bullet: Método / Critério / Norma de recuperação),
# recomputing character offsets fresh (since earlier edits may have shifted them).
function Get-LabeledRunText($label, $nextLabel) {
    $p = $d.Paragraphs.Item(17)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End

    $rLabel = $d.Range($pStart, $pEnd)
    [void]$rLabel.Find.Execute($label, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $textStart = $rLabel.End

    $textEnd = $pEnd
    if ($nextLabel) {
        $rNext = $d.Range($textStart, $pEnd)
        [void]$rNext.Find.Execute($nextLabel, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        $textEnd = $rNext.Start
    }

    return $d.Range($textStart, $textEnd).Text
}

# Sets the text of the run-range between two labels inside paragraph #17,
# recomputing character offsets fresh right before the write.
function Set-LabeledRunText($label, $nextLabel, $newText) {
    $p = $d.Paragraphs.Item(17)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End

    $rLabel = $d.Range($pStart, $pEnd)
    [void]$rLabel.Find.Execute($label, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $textStart = $rLabel.End

    $textEnd = $pEnd
    if ($nextLabel) {
        $rNext = $d.Range($textStart, $pEnd)
        [void]$rNext.Find.Execute($nextLabel, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        $textEnd = $rNext.Start
    }

    $target = $d.Range($textStart, $textEnd)
    $target.Text = $newText
}

# --- Capture all "before" values we need, before making any changes ---

$objetivosPT        = Get-ParaText 6    # "A disciplina tem o objetivo..."
$objetivosEN        = Get-ParaText 7    # "The course aims to introduce..."
$docenteNome        = Get-ParaText 9    # "5840942 - Marco Aurélio..."
$programaResumidoPT = Get-ParaText 11   # "Introdução. Formação do solo..."
$programaResumidoEN = Get-ParaText 12   # "Introduction. Soil formation..."
$programaPT         = Get-ParaText 14   # "INTRODUÇÃO. Conceitos Básicos..."
$bibliografiaBody   = Get-ParaText 19   # "Bibliografia básica:1. MEURER..."

$metodoTextRaw   = Get-LabeledRunText "Método: " "Critério: "                  # "O aluno poderá optar..." (+ trailing line-break char)
$criterioText    = Get-LabeledRunText "Critério: " "Norma de recuperação: "    # "Exame Final (EF)..." (+ trailing line-break char)

# The captured run text above includes a trailing vertical-tab char (chr 11)
# representing the <w:br/> that ends those runs. When moving this text into a
# plain paragraph (#14) that has no line break, that trailing char must be stripped.
$metodoText = $metodoTextRaw
if ($metodoText.Length -gt 0 -and [int][char]$metodoText.Substring($metodoText.Length - 1) -eq 11) {
    $metodoText = $metodoText.Substring(0, $metodoText.Length - 1)
}

# --- Apply the rotation: each slot receives the value that, per the diff, now belongs there ---

$d.Paragraphs.Item(6).Range.Text = $programaResumidoPT     # Objetivos PT         -> Programa resumido PT text
$d.Paragraphs.Item(7).Range.Text = $programaResumidoEN     # Objetivos EN         -> Programa resumido EN text
$d.Paragraphs.Item(9).Range.Text = $objetivosPT             # Docente name         -> old Objetivos PT text
$d.Paragraphs.Item(11).Range.Text = $programaPT             # Programa resumido PT -> old Programa PT text
$d.Paragraphs.Item(12).Range.Text = $objetivosEN            # Programa resumido EN -> old Objetivos EN text
$d.Paragraphs.Item(14).Range.Text = $metodoText              # Programa PT          -> old Método text
$d.Paragraphs.Item(19).Range.Text = $docenteNome             # Bibliografia body    -> old Docente name text

# Sub-runs of paragraph 17 ("Avaliação"): Método gets old Critério text,
# Norma de recuperação gets old Bibliografia body text.
Set-LabeledRunText "Método: " "Critério: " $criterioText
Set-LabeledRunText "Norma de recuperação: " $null $bibliografiaBody
